$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.035.21'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.797.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '662.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.45'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.796.12'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.52%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.52%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.99'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.91%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.63'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.438.08'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.802.76'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.990.19'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.73'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '473.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.19'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.712'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.86%  '

$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000145'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.40%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.65'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.22'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.97%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.56%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.948.14'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.34%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.93%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.35'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.07'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.178'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +18.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.753.24'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.43%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.77%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.31'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.88%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.16%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.07'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '159.04'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.90'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.16%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.43'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.07%  '

$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.301'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.48'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.30%  '
